# Add 2022-Q4 data:
#  - insert a new "2022-Q4" worksheet (fund-holdings table) before "2022-Q1"
#  - update the "总计" summary sheet: row 2 becomes the new 2022-Q4 summary
#    and the old 2022-Q1 summary row is pushed down to row 3

$wb = $excel.ActiveWorkbook

$fmtFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------------
# Helper-ish references: style sources already present in the workbook.
#   Range("A2") on "总计"  -> cellXf s="2" applied to a NUMBER cell
#   Range("B1") on "总计"  -> cellXf s="2" applied to a TEXT (header) cell
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$numStyleSrc  = $summary.Range("A2")
$textStyleSrc = $summary.Range("B1")

function Set-TextValue {
    param($range, [string]$text)
    # Force literal text storage (so numeric-looking strings like "010583"
    # or "13.38" keep leading zeros / stay text instead of becoming numbers),
    # then drop back to the default "Normal" style so no stray number format
    # lingers on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

function Set-StyledText {
    param($range, [string]$text, $styleSrc)
    Set-TextValue $range $text
    $styleSrc.Copy()
    $range.PasteSpecial($fmtFormats)
}

function Set-StyledNumber {
    param($range, $num, $styleSrc)
    $range.Value = $num
    $styleSrc.Copy()
    $range.PasteSpecial($fmtFormats)
}

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right before the existing "2022-Q1".
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q4 = $wb.Worksheets.Add($q1)
$q4.Name = "2022-Q4"

# Header row (all header cells share the "总计"!B1 style -> s="2")
Set-StyledText $q4.Range("B1") "基金代码"         $textStyleSrc
Set-StyledText $q4.Range("C1") "基金名称"         $textStyleSrc
Set-StyledText $q4.Range("D1") "基金规模"         $textStyleSrc
Set-StyledText $q4.Range("E1") "股票总仓位"       $textStyleSrc
Set-StyledText $q4.Range("F1") "仓位占比"         $textStyleSrc
Set-StyledText $q4.Range("G1") "持有市值(亿元)"   $textStyleSrc
Set-StyledText $q4.Range("H1") "仓位排名"         $textStyleSrc

# Data rows. Column A carries the "总计"!A2 style -> s="2"; the rest are
# plain (default style). Columns B-G are text (even the numeric-looking
# ones), column H is a real number.
Set-StyledNumber $q4.Range("A2") 0 $numStyleSrc
Set-TextValue    $q4.Range("B2") "010583"
Set-TextValue    $q4.Range("C2") "富国蓝筹精选股票（QDII）美元"
Set-TextValue    $q4.Range("D2") "13.38"
Set-TextValue    $q4.Range("E2") "91.97"
Set-TextValue    $q4.Range("F2") "4.91"
Set-TextValue    $q4.Range("G2") "0.6570"
$q4.Range("H2").Value = 3

Set-StyledNumber $q4.Range("A3") 1 $numStyleSrc
Set-TextValue    $q4.Range("B3") "007455"
Set-TextValue    $q4.Range("C3") "富国蓝筹精选股票（QDII）人民币"
Set-TextValue    $q4.Range("D3") "13.38"
Set-TextValue    $q4.Range("E3") "91.97"
Set-TextValue    $q4.Range("F3") "4.91"
Set-TextValue    $q4.Range("G3") "0.6570"
$q4.Range("H3").Value = 3

Set-StyledNumber $q4.Range("A4") 2 $numStyleSrc
Set-TextValue    $q4.Range("B4") "100055"
Set-TextValue    $q4.Range("C4") "富国全球科技互联网股票（QDII）"
Set-TextValue    $q4.Range("D4") "3.86"
Set-TextValue    $q4.Range("E4") "94.32"
Set-TextValue    $q4.Range("F4") "4.90"
Set-TextValue    $q4.Range("G4") "0.1891"
$q4.Range("H4").Value = 6

Set-StyledNumber $q4.Range("A5") 3 $numStyleSrc
Set-TextValue    $q4.Range("B5") "006781"
Set-TextValue    $q4.Range("C5") "汇丰晋信港股通精选股票"
Set-TextValue    $q4.Range("D5") "0.70"
Set-TextValue    $q4.Range("E5") "91.29"
Set-TextValue    $q4.Range("F5") "2.98"
Set-TextValue    $q4.Range("G5") "0.0209"
$q4.Range("H5").Value = 10

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet.
#    Push the existing 2022-Q1 summary row down to row 3, then overwrite
#    row 2 with the new 2022-Q4 summary.
# ---------------------------------------------------------------------------
Set-StyledNumber $summary.Range("A3") 1 $numStyleSrc
Set-TextValue    $summary.Range("B3") "2022-Q1"
$summary.Range("C3").Value = 8
$summary.Range("D3").Value = 1.87

Set-TextValue $summary.Range("B2") "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 1.52

# Keep "总计" as the active/selected sheet (matches the unchanged bookViews
# activeTab="0" from the original workbook).
$summary.Select()
